$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D column values are numeric-looking strings (e.g. "1.005", "10.00") that must
# stay plain text exactly as scraped. Setting NumberFormat to Text ("@") before
# assigning the literal, then resetting the style back to Normal afterwards,
# keeps the stored string intact without leaving a residual text-format style on
# the cell (matching the original workbook, where these cells carry no style).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.865.79"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.27%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.881.84"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.31%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.14%  "
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4673"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.80%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3937"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07931"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9772"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.83%  "
$ws.Range("E11").Value = "  +2.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.881.68"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.748"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.88%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.008"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06973"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.67"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.31%  "
$ws.Range("E17").Value = "  +0.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001011"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.64%  "
$ws.Range("E20").Value = "  +0.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "28.873.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.365"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.91%  "
$ws.Range("E23").Value = "  +1.71%  "
$ws.Range("E24").Value = "  +1.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.136.82"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.58%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.770"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.27%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.008"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.95%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09401"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9417"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.314"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.354"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.351"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05924"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.37%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02124"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.149"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.932"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5726"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.11%  "
$ws.Range("E42").Value = "  +1.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.07257"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.89"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.78%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5351"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.150"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.128"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.853"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.64%  "
$ws.Range("E49").Value = "  +1.67%  "
$ws.Range("E50").Value = "  +3.33%  "
$ws.Range("E51").Value = "  +0.45%  "
